$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*September 19, 2025*") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# ---------------------------------------------------------------------
# 2. Split the single "street, city state zip" address line into two
#    separate lines:
#        "2933 Lamory Pl"
#        "Santa Clara, CA 95051"
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*2933 Lamory Pl, Santa Clara CA 95051*") {
        # Shrink the existing paragraph down to just the street address.
        $p.Range.Text = "2933 Lamory Pl"

        # Create a brand-new paragraph right after it (inherits the same
        # paragraph formatting) to hold the city/state/zip line.
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()

        # Insert the paragraph's contents via raw WordprocessingML so the
        # xml:space="preserve" attribute on <w:t> is retained exactly like
        # the rest of the document's runs.
        $paraXml = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Santa Clara, CA 95051</w:t></w:r></w:p>'
        $newPara.Range.InsertXML($paraXml)
        break
    }
}

# ---------------------------------------------------------------------
# 3. Remove the now-superfluous blank "NoSpacing" paragraph that sits
#    directly after the "...Board of Directors" signature line.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Board of Directors*") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text.Trim() -eq "") {
            $next.Range.Delete()
        }
        break
    }
}
